$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 950.6875
$ws.Cells.Item(41, 9).Value = 846.4666999999999
$ws.Cells.Item(41, 10).Value = 1042.6471
$ws.Cells.Item(41, 11).Value = 846.4666999999999
$ws.Cells.Item(41, 12).Value = 1042.6471
$ws.Cells.Item(41, 13).Value = -406.4666999999999
$ws.Cells.Item(41, 14).Value = -1922.6471

$ws.Cells.Item(55, 8).Value = 182.28572
$ws.Cells.Item(55, 9).Value = 197.33333
$ws.Cells.Item(55, 10).Value = 92
$ws.Cells.Item(55, 11).Value = 197.33333
$ws.Cells.Item(55, 12).Value = 92
$ws.Cells.Item(55, 13).Value = 16.66667000000001
$ws.Cells.Item(55, 14).Value = -520

$ws.Cells.Item(70, 8).Value = 9744.923000000001
$ws.Cells.Item(70, 10).Value = 11819.3
$ws.Cells.Item(70, 12).Value = 35457.89999999999
$ws.Cells.Item(70, 14).Value = -35997.89999999999

$ws.Cells.Item(73, 8).Value = 9744.923000000001
$ws.Cells.Item(73, 10).Value = 11819.3
$ws.Cells.Item(73, 12).Value = 35457.89999999999
$ws.Cells.Item(73, 14).Value = -37329.89999999999

$ws.Cells.Item(74, 8).Value = 6639.8887
$ws.Cells.Item(74, 10).Value = 6685.7
$ws.Cells.Item(74, 12).Value = 6685.7
$ws.Cells.Item(74, 14).Value = -8557.700000000001

$ws.Cells.Item(76, 8).Value = 77696140
$ws.Cells.Item(76, 9).Value = 126253150
$ws.Cells.Item(76, 11).Value = 126253150
$ws.Cells.Item(76, 13).Value = -126252835

$ws.Cells.Item(77, 8).Value = 6639.8887
$ws.Cells.Item(77, 10).Value = 6685.7
$ws.Cells.Item(77, 12).Value = 33428.5
$ws.Cells.Item(77, 14).Value = -42788.5

$ws.Cells.Item(79, 8).Value = 77696140
$ws.Cells.Item(79, 9).Value = 126253150
$ws.Cells.Item(79, 11).Value = 126253150
$ws.Cells.Item(79, 13).Value = -126252058

$ws.Cells.Item(96, 8).Value = 864.4286
$ws.Cells.Item(96, 10).Value = 923.25
$ws.Cells.Item(96, 12).Value = 2769.75
$ws.Cells.Item(96, 14).Value = -5515.75

$ws.Cells.Item(112, 8).Value = 2685.6667
$ws.Cells.Item(112, 10).Value = 2720.818
$ws.Cells.Item(112, 12).Value = 8162.454000000001
$ws.Cells.Item(112, 14).Value = -10378.454

$ws.Cells.Item(132, 8).Value = 3772.9778
$ws.Cells.Item(132, 10).Value = 2570.7856
$ws.Cells.Item(132, 12).Value = 7712.3568
$ws.Cells.Item(132, 14).Value = -12772.3568

$ws.Cells.Item(137, 8).Value = 4712
$ws.Cells.Item(137, 9).Value = 3796.8
$ws.Cells.Item(137, 11).Value = 11390.4
$ws.Cells.Item(137, 13).Value = -8840.400000000001

$ws.Cells.Item(138, 8).Value = 4747.375
$ws.Cells.Item(138, 10).Value = 5331.273
$ws.Cells.Item(138, 12).Value = 15993.819
$ws.Cells.Item(138, 14).Value = -26273.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2574.48
$ws.Cells.Item(45, 9).Value = 1618.4667
$ws.Cells.Item(45, 11).Value = 1618.4667
$ws.Cells.Item(45, 13).Value = -1241.4667

$ws.Cells.Item(74, 8).Value = 2142
$ws.Cells.Item(74, 9).Value = 2142
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 2142
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -1268
$ws.Cells.Item(74, 14).ClearContents()

$ws.Cells.Item(77, 8).Value = 2142
$ws.Cells.Item(77, 9).Value = 2142
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 10710
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -6342
$ws.Cells.Item(77, 14).ClearContents()

$ws.Cells.Item(88, 8).Value = 35542
$ws.Cells.Item(88, 9).Value = 36524.5
$ws.Cells.Item(88, 10).Value = 35050.75
$ws.Cells.Item(88, 11).Value = 36524.5
$ws.Cells.Item(88, 12).Value = 35050.75
$ws.Cells.Item(88, 13).Value = -36118.5
$ws.Cells.Item(88, 14).Value = -35862.75

$ws.Cells.Item(91, 8).Value = 35542
$ws.Cells.Item(91, 9).Value = 36524.5
$ws.Cells.Item(91, 10).Value = 35050.75
$ws.Cells.Item(91, 11).Value = 36524.5
$ws.Cells.Item(91, 12).Value = 35050.75
$ws.Cells.Item(91, 13).Value = -35120.5
$ws.Cells.Item(91, 14).Value = -37858.75

$ws.Cells.Item(102, 8).Value = 1061.963
$ws.Cells.Item(102, 9).Value = 1045.1154
$ws.Cells.Item(102, 11).Value = 1045.1154
$ws.Cells.Item(102, 13).Value = 576.8846000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3152
$ws.Cells.Item(20, 9).Value = 1457
$ws.Cells.Item(20, 10).Value = 3999.5
$ws.Cells.Item(20, 11).Value = 1457
$ws.Cells.Item(20, 12).Value = 3999.5
$ws.Cells.Item(20, 13).Value = -1210
$ws.Cells.Item(20, 14).Value = -4493.5

$ws.Cells.Item(26, 8).Value = 41595.6
$ws.Cells.Item(26, 9).Value = 41595.6
$ws.Cells.Item(26, 11).Value = 41595.6
$ws.Cells.Item(26, 13).Value = -41303.6

$ws.Cells.Item(86, 8).Value = 26939.889
$ws.Cells.Item(86, 9).Value = 7797.1113
$ws.Cells.Item(86, 11).Value = 7797.1113
$ws.Cells.Item(86, 13).Value = -6674.1113

$ws.Cells.Item(89, 8).Value = 26939.889
$ws.Cells.Item(89, 9).Value = 7797.1113
$ws.Cells.Item(89, 11).Value = 38985.5565
$ws.Cells.Item(89, 13).Value = -33369.5565

$ws.Cells.Item(99, 8).Value = 2448
$ws.Cells.Item(99, 9).Value = 2448
$ws.Cells.Item(99, 11).Value = 2448
$ws.Cells.Item(99, 13).Value = -950

$ws.Cells.Item(105, 8).Value = 3632.7778
$ws.Cells.Item(105, 9).Value = 3340.6
$ws.Cells.Item(105, 10).Value = 3998
$ws.Cells.Item(105, 11).Value = 3340.6
$ws.Cells.Item(105, 12).Value = 3998
$ws.Cells.Item(105, 13).Value = -1593.6
$ws.Cells.Item(105, 14).Value = -7492

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 80511
$ws.Cells.Item(52, 10).Value = 80511
$ws.Cells.Item(52, 12).Value = 80511
$ws.Cells.Item(52, 14).Value = -81099

$ws.Cells.Item(62, 8).Value = 22031.666
$ws.Cells.Item(62, 9).Value = 5787.25
$ws.Cells.Item(62, 11).Value = 5787.25
$ws.Cells.Item(62, 13).Value = -5163.25

$ws.Cells.Item(65, 8).Value = 22031.666
$ws.Cells.Item(65, 9).Value = 5787.25
$ws.Cells.Item(65, 11).Value = 28936.25
$ws.Cells.Item(65, 13).Value = -25816.25

$ws.Cells.Item(135, 8).Value = 85999.2
$ws.Cells.Item(135, 9).Value = 9999
$ws.Cells.Item(135, 10).Value = 104999.25
$ws.Cells.Item(135, 11).Value = 9999
$ws.Cells.Item(135, 12).Value = 104999.25
$ws.Cells.Item(135, 13).Value = -4929
$ws.Cells.Item(135, 14).Value = -115139.25

$ws.Cells.Item(138, 8).Value = 63737.09
$ws.Cells.Item(138, 10).Value = 67499.75
$ws.Cells.Item(138, 12).Value = 67499.75
$ws.Cells.Item(138, 14).Value = -77779.75

$ws.Cells.Item(140, 8).Value = 139999
$ws.Cells.Item(140, 10).Value = 139999
$ws.Cells.Item(140, 12).Value = 139999
$ws.Cells.Item(140, 14).Value = -150359

$ws.Cells.Item(141, 8).Value = 58808
$ws.Cells.Item(141, 10).Value = 58808
$ws.Cells.Item(141, 12).Value = 58808
$ws.Cells.Item(141, 14).Value = -69168

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 3515.8823
$ws.Cells.Item(86, 9).Value = 2745.5
$ws.Cells.Item(86, 11).Value = 8236.5
$ws.Cells.Item(86, 13).Value = -7050.5

$ws.Cells.Item(89, 8).Value = 3515.8823
$ws.Cells.Item(89, 9).Value = 2745.5
$ws.Cells.Item(89, 11).Value = 24709.5
$ws.Cells.Item(89, 13).Value = -18781.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8253.166999999999
$ws.Cells.Item(70, 9).Value = 7979.9165
$ws.Cells.Item(70, 10).Value = 8799.666999999999
$ws.Cells.Item(70, 11).Value = 7979.9165
$ws.Cells.Item(70, 12).Value = 8799.666999999999
$ws.Cells.Item(70, 13).Value = -7709.9165
$ws.Cells.Item(70, 14).Value = -9339.666999999999

$ws.Cells.Item(73, 8).Value = 8253.166999999999
$ws.Cells.Item(73, 9).Value = 7979.9165
$ws.Cells.Item(73, 10).Value = 8799.666999999999
$ws.Cells.Item(73, 11).Value = 7979.9165
$ws.Cells.Item(73, 12).Value = 8799.666999999999
$ws.Cells.Item(73, 13).Value = -7043.9165
$ws.Cells.Item(73, 14).Value = -10671.667

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2284.1875
$ws.Cells.Item(46, 9).Value = 1700.2222
$ws.Cells.Item(46, 10).Value = 3035
$ws.Cells.Item(46, 11).Value = 1700.2222
$ws.Cells.Item(46, 12).Value = 3035
$ws.Cells.Item(46, 13).Value = -1512.2222
$ws.Cells.Item(46, 14).Value = -3411

$ws.Cells.Item(61, 8).Value = 2476.9
$ws.Cells.Item(61, 9).Value = 2441.111
$ws.Cells.Item(61, 10).Value = 2799
$ws.Cells.Item(61, 11).Value = 2441.111
$ws.Cells.Item(61, 12).Value = 2799
$ws.Cells.Item(61, 13).Value = -2239.111
$ws.Cells.Item(61, 14).Value = -3203

$ws.Cells.Item(93, 8).Value = 2103.2856
$ws.Cells.Item(93, 9).Value = 1885.6666
$ws.Cells.Item(93, 10).Value = 2495
$ws.Cells.Item(93, 11).Value = 1885.6666
$ws.Cells.Item(93, 12).Value = 2495
$ws.Cells.Item(93, 13).Value = -637.6666
$ws.Cells.Item(93, 14).Value = -4991

$ws.Cells.Item(113, 8).Value = 2476.9
$ws.Cells.Item(113, 9).Value = 2441.111
$ws.Cells.Item(113, 10).Value = 2799
$ws.Cells.Item(113, 11).Value = 2441.111
$ws.Cells.Item(113, 12).Value = 2799
$ws.Cells.Item(113, 13).Value = -271.1109999999999
$ws.Cells.Item(113, 14).Value = -7139

$ws.Cells.Item(132, 8).Value = 2683.0312
$ws.Cells.Item(132, 9).Value = 1780.3889
$ws.Cells.Item(132, 10).Value = 3843.5715
$ws.Cells.Item(132, 11).Value = 5341.1667
$ws.Cells.Item(132, 12).Value = 11530.7145
$ws.Cells.Item(132, 13).Value = -2811.1667
$ws.Cells.Item(132, 14).Value = -16590.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 4813.5713
$ws.Cells.Item(96, 9).Value = 4769
$ws.Cells.Item(96, 10).Value = 4821
$ws.Cells.Item(96, 11).Value = 4769
$ws.Cells.Item(96, 12).Value = 4821
$ws.Cells.Item(96, 13).Value = -3396
